$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.369.82'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.252.94'
$ws.Range("E3").Value = '  +3.84%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.47'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.06'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.248.34'
$ws.Range("E8").Value = '  +3.96%  '
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("E13").Value = '  -2.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.22'
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.783.65'
$ws.Range("E15").Value = '  +4.06%  '
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.247.66'
$ws.Range("E17").Value = '  +3.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.416.91'
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '474.51'
$ws.Range("E20").Value = '  -1.65%  '
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("E22").Value = '  +3.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.95'
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.89'
$ws.Range("E24").Value = '  -4.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.13'
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  +3.87%  '
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("E30").Value = '  +3.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.58'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("E33").Value = '  -4.08%  '
$ws.Range("E34").Value = '  -2.99%  '
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("E36").Value = '  -1.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.73'
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0708'
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '421.84'
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.983.78'
$ws.Range("E41").Value = '  +2.72%  '
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("E43").Value = '  -5.17%  '
$ws.Range("E44").Value = '  -7.33%  '
$ws.Range("E45").Value = '  +2.47%  '
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.77'
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.30'
$ws.Range("E49").Value = '  -2.43%  '
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '122.31'
$ws.Range("E51").Value = '  +1.54%  '
